# Auto-generated edit script for QualificationPAC workbook update
$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Date and Base Definition ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-07-21T11:52:46+00:00"
$meta.Range("B18").Value = "https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/SavoirFaire"

# --- Sheet "Elements": update row 3 (was qualificationPAC, becomes typeSavoirFaire) ---
$ws = $wb.Worksheets.Item("Elements")
$ws.Range("A3").Value = "QualificationPAC.typeSavoirFaire"
$ws.Range("B3").Value = "QualificationPAC.typeSavoirFaire"
$ws.Range("L3").Value = " Le type de savoir-faire (qualifications/autres attributions) désigne par exemple:** une spécialité ordinale (S);** une compétence (C);** etc."
$ws.Range("M3").Value = " Le type de savoir-faire (qualifications/autres attributions) désigne par exemple:** une spécialité ordinale (S);** une compétence (C);** etc."
$ws.Range("Z3").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R04-TypeSavoirFaire?vs"
$ws.Range("AF3").Value = "SavoirFaire.typeSavoirFaire"

# --- Add new rows 4, 5, 6; copy formatting from row 3 first ---
$ws.Range("A3:AJ3").Copy()
$ws.Range("A4:AJ6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 4: QualificationPAC.dateReconnaissance
$ws.Range("A4").Value = "QualificationPAC.dateReconnaissance"
$ws.Range("B4").Value = "QualificationPAC.dateReconnaissance"
$ws.Range("D4").Value = ""
$ws.Range("F4").Value = "0"
$ws.Range("G4").Value = "1"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = "date`n"
$ws.Range("L4").Value = " Date à laquelle, l’organisme donnant l’autorisation d’exercer une qualification a reconnu cette qualification ou date à laquelle l'attribution a été donnée au professionnel."
$ws.Range("M4").Value = " Date à laquelle, l’organisme donnant l’autorisation d’exercer une qualification a reconnu cette qualification ou date à laquelle l'attribution a été donnée au professionnel."
$ws.Range("P4").Value = ""
$ws.Range("R4").Value = ""
$ws.Range("S4").Value = ""
$ws.Range("T4").Value = ""
$ws.Range("U4").Value = ""
$ws.Range("V4").Value = ""
$ws.Range("W4").Value = ""
$ws.Range("X4").Value = ""
$ws.Range("Y4").Value = ""
$ws.Range("Z4").Value = ""
$ws.Range("AA4").Value = ""
$ws.Range("AB4").Value = ""
$ws.Range("AC4").Value = ""
$ws.Range("AD4").Value = ""
$ws.Range("AE4").Value = ""
$ws.Range("AF4").Value = "SavoirFaire.dateReconnaissance"
$ws.Range("AG4").Value = "0"
$ws.Range("AH4").Value = "1"
$ws.Range("AI4").Value = ""
$ws.Range("AJ4").Value = ""

# Row 5: QualificationPAC.dateAbandon
$ws.Range("A5").Value = "QualificationPAC.dateAbandon"
$ws.Range("B5").Value = "QualificationPAC.dateAbandon"
$ws.Range("D5").Value = ""
$ws.Range("F5").Value = "0"
$ws.Range("G5").Value = "1"
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = "date`n"
$ws.Range("L5").Value = " Date à laquelle le professionnel a déclaré renoncer à l’exercice d’un savoir-faire ou date à laquelle il ne souhaite plus le faire apparaître."
$ws.Range("M5").Value = " Date à laquelle le professionnel a déclaré renoncer à l’exercice d’un savoir-faire ou date à laquelle il ne souhaite plus le faire apparaître."
$ws.Range("P5").Value = ""
$ws.Range("R5").Value = ""
$ws.Range("S5").Value = ""
$ws.Range("T5").Value = ""
$ws.Range("U5").Value = ""
$ws.Range("V5").Value = ""
$ws.Range("W5").Value = ""
$ws.Range("X5").Value = ""
$ws.Range("Y5").Value = ""
$ws.Range("Z5").Value = ""
$ws.Range("AA5").Value = ""
$ws.Range("AB5").Value = ""
$ws.Range("AC5").Value = ""
$ws.Range("AD5").Value = ""
$ws.Range("AE5").Value = ""
$ws.Range("AF5").Value = "SavoirFaire.dateAbandon"
$ws.Range("AG5").Value = "0"
$ws.Range("AH5").Value = "1"
$ws.Range("AI5").Value = ""
$ws.Range("AJ5").Value = ""

# Row 6: QualificationPAC.qualificationPAC (moved down from old row 3)
$ws.Range("A6").Value = "QualificationPAC.qualificationPAC"
$ws.Range("B6").Value = "QualificationPAC.qualificationPAC"
$ws.Range("D6").Value = ""
$ws.Range("F6").Value = "0"
$ws.Range("G6").Value = "1"
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = "Coding`n"
$ws.Range("L6").Value = " Qualification de praticien adjoint contractuel."
$ws.Range("M6").Value = " Qualification de praticien adjoint contractuel."
$ws.Range("P6").Value = ""
$ws.Range("R6").Value = ""
$ws.Range("S6").Value = ""
$ws.Range("T6").Value = ""
$ws.Range("U6").Value = ""
$ws.Range("V6").Value = ""
$ws.Range("W6").Value = ""
$ws.Range("X6").Value = "preferred"
$ws.Range("Z6").Value = "https://interop.esante.gouv.fr/ig/fhir/mos/ValueSet/qualificationPAC-vs"
$ws.Range("AA6").Value = ""
$ws.Range("AB6").Value = ""
$ws.Range("AC6").Value = ""
$ws.Range("AD6").Value = ""
$ws.Range("AE6").Value = ""
$ws.Range("AF6").Value = "QualificationPAC.qualificationPAC"
$ws.Range("AG6").Value = "0"
$ws.Range("AH6").Value = "1"
$ws.Range("AI6").Value = ""
$ws.Range("AJ6").Value = ""

# --- Column width adjustments (best-effort, quantized by Excel column-width units) ---
$ws.Columns.Item(1).ColumnWidth = 29.833333333333332
$ws.Columns.Item(2).ColumnWidth = 29.833333333333332
$ws.Columns.Item(26).ColumnWidth = 68.0

"done"
